$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preventing Excel from
# auto-converting number-looking strings (e.g. "1.12") into real numbers.
# We borrow the (always-plain) style from column E of the same row, which
# keeps the cell on style index 0 -- i.e. no visible / structural style change.
function Set-TextValue {
    param($cell, $value, $styleSource)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $styleSource.Style
}

# Row 2
$ws.Range("D2").Value = "91.643.96"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "3.141.10"
$ws.Range("E3").Value = "  +1.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "241.01" $ws.Range("E5")
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
Set-TextValue $ws.Range("D6") "619.28" $ws.Range("E6")
$ws.Range("E6").Value = "  -0.89%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.12" $ws.Range("E7")
$ws.Range("E7").Value = "  -2.95%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.389" $ws.Range("E8")
$ws.Range("E8").Value = "  +4.79%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").Value = "3.136.83"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.750" $ws.Range("E11")
$ws.Range("E11").Value = "  +0.30%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.204" $ws.Range("E12")
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000255" $ws.Range("E13")
$ws.Range("E13").Value = "  +1.27%  "

# Row 14
Set-TextValue $ws.Range("D14") "35.16" $ws.Range("E14")
$ws.Range("E14").Value = "  -0.23%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.62" $ws.Range("E15")
$ws.Range("E15").Value = "  +2.65%  "

# Row 16
$ws.Range("D16").Value = "91.297.84"
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").Value = "3.715.97"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18
$ws.Range("D18").Value = "3.160.55"
$ws.Range("E18").Value = "  +1.75%  "

# Row 19
Set-TextValue $ws.Range("D19") "3.77" $ws.Range("E19")
$ws.Range("E19").Value = "  -1.39%  "

# Row 20
Set-TextValue $ws.Range("D20") "14.96" $ws.Range("E20")
$ws.Range("E20").Value = "  +4.78%  "

# Row 21
Set-TextValue $ws.Range("D21") "5.91" $ws.Range("E21")
$ws.Range("E21").Value = "  +2.30%  "

# Row 22
Set-TextValue $ws.Range("D22") "457.31" $ws.Range("E22")
$ws.Range("E22").Value = "  +2.82%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.0000202" $ws.Range("E23")
$ws.Range("E23").Value = "  -3.46%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.21" $ws.Range("E24")
$ws.Range("E24").Value = "  +1.74%  "

# Row 25
Set-TextValue $ws.Range("D25") "5.93" $ws.Range("E25")
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
Set-TextValue $ws.Range("D26") "88.55" $ws.Range("E26")
$ws.Range("E26").Value = "  -4.55%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.87" $ws.Range("E27")
$ws.Range("E27").Value = "  -0.99%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.152" $ws.Range("E28")
$ws.Range("E28").Value = "  +42.63%  "

# Row 29
$ws.Range("E29").Value = "  +1.70%  "

# Row 30
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.231" $ws.Range("E31")
$ws.Range("E31").Value = "  +3.94%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.167" $ws.Range("E32")
$ws.Range("E32").Value = "  -6.18%  "

# Row 33
Set-TextValue $ws.Range("D33") "9.36" $ws.Range("E33")
$ws.Range("E33").Value = "  +1.45%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.176" $ws.Range("E34")
$ws.Range("E34").Value = "  +11.61%  "

# Row 35
Set-TextValue $ws.Range("D35") "26.41" $ws.Range("E35")
$ws.Range("E35").Value = "  -0.48%  "

# Row 36
Set-TextValue $ws.Range("D36") "7.44" $ws.Range("E36")
$ws.Range("E36").Value = "  -2.44%  "

# Row 37
$ws.Range("E37").Value = "  +4.29%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D38") "492.93" $ws.Range("E38")
$ws.Range("E38").Value = "  +0.35%  "

# Row 39
$ws.Range("B39").Value = "MantraDAO"
$ws.Range("C39").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D39") "3.90" $ws.Range("E39")
$ws.Range("E39").Value = "  -9.39%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.32" $ws.Range("E40")
$ws.Range("E40").Value = "  +2.50%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.441" $ws.Range("E41")
$ws.Range("E41").Value = "  +6.07%  "

# Row 42
Set-TextValue $ws.Range("D42") "3.40" $ws.Range("E42")
$ws.Range("E42").Value = "  -5.08%  "

# Row 43
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D43") "22.14" $ws.Range("E43")
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("B44").Value = "Binance-PegBSC-USD"
$ws.Range("C44").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D44") "0.751" $ws.Range("E44")
$ws.Range("E44").Value = "  -24.78%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D46") "1.93" $ws.Range("E46")
$ws.Range("E46").Value = "  +1.93%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D47") "0.707" $ws.Range("E47")
$ws.Range("E47").Value = "  +3.62%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D48") "156.77" $ws.Range("E48")
$ws.Range("E48").Value = "  -1.37%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.36" $ws.Range("E49")
$ws.Range("E49").Value = "  +2.12%  "

# Row 50
Set-TextValue $ws.Range("D50") "4.48" $ws.Range("E50")
$ws.Range("E50").Value = "  -1.65%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D51") "0.0327" $ws.Range("E51")
$ws.Range("E51").Value = "  +5.33%  "
